# Add a new "Compact List" paragraph style (styleId "CompactList"), cloned
# from the existing "Compact" style: based on Body Text, quick-style, with
# 36-twip (1.8pt) spacing before/after.

$d = $word.ActiveDocument

$wdStyleTypeParagraph = 1

$newStyle = $d.Styles.Add("CompactList", $wdStyleTypeParagraph)
$newStyle.NameLocal = "Compact List"
$newStyle.BaseStyle = "BodyText"
$newStyle.QuickStyle = $true
$newStyle.ParagraphFormat.SpaceBefore = 1.8
$newStyle.ParagraphFormat.SpaceAfter = 1.8
